# Auto-generated Excel COM-interop script to apply diff changes
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 76900
$ws.Range("J68").Value = 76900
$ws.Range("L68").Value = 76900
$ws.Range("N68").Value = -78398
$ws.Range("H71").Value = 76900
$ws.Range("J71").Value = 76900
$ws.Range("L71").Value = 230700
$ws.Range("N71").Value = -238188
$ws.Range("H96").Value = 17363
$ws.Range("I96").Value = 1060.5
$ws.Range("J96").Value = 39099.668
$ws.Range("K96").Value = 3181.5
$ws.Range("L96").Value = 117299.004
$ws.Range("M96").Value = -1808.5
$ws.Range("N96").Value = -120045.004
$ws.Range("H113").Value = 1674.5
$ws.Range("I113").Value = 1566
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1566
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1688
$ws.Range("N113").Value = -8508
$ws.Range("H116").Value = 41690220
$ws.Range("J116").Value = 7219.5
$ws.Range("L116").Value = 7219.5
$ws.Range("N116").Value = -14103.5
$ws.Range("H137").Value = 3860
$ws.Range("I137").Value = 3109.5833
$ws.Range("K137").Value = 9328.749899999999
$ws.Range("M137").Value = -6778.749899999999
$ws.Range("H138").Value = 6311.4053
$ws.Range("J138").Value = 6855.778
$ws.Range("L138").Value = 20567.334
$ws.Range("N138").Value = -30847.334

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2871832.8
$ws.Range("I32").Value = 3854930.5
$ws.Range("J32").Value = 31773.223
$ws.Range("K32").Value = 3854930.5
$ws.Range("L32").Value = 31773.223
$ws.Range("M32").Value = -3854643.5
$ws.Range("N32").Value = -32347.223
$ws.Range("H61").Value = 5474.4287
$ws.Range("I61").Value = 3831.2
$ws.Range("K61").Value = 3831.2
$ws.Range("M61").Value = -3619.2
$ws.Range("H74").Value = 3476.9285
$ws.Range("I74").Value = 2974.16
$ws.Range("K74").Value = 2974.16
$ws.Range("M74").Value = -2100.16
$ws.Range("H77").Value = 3476.9285
$ws.Range("I77").Value = 2974.16
$ws.Range("K77").Value = 14870.8
$ws.Range("M77").Value = -10502.8
$ws.Range("H118").Value = 55500
$ws.Range("J118").Value = 55500
$ws.Range("L118").Value = 55500
$ws.Range("N118").Value = -58814
$ws.Range("H122").Value = 286527.84
$ws.Range("I122").Value = 448460.75
$ws.Range("J122").Value = 8928.571
$ws.Range("K122").Value = 1345382.25
$ws.Range("L122").Value = 26785.713
$ws.Range("M122").Value = -1342932.25
$ws.Range("N122").Value = -31685.713
$ws.Range("H132").Value = 230155.48
$ws.Range("I132").Value = 402293.4
$ws.Range("J132").Value = 3658.2104
$ws.Range("K132").Value = 1206880.2
$ws.Range("L132").Value = 10974.6312
$ws.Range("M132").Value = -1204350.2
$ws.Range("N132").Value = -16034.6312
$ws.Range("H136").Value = 5474.4287
$ws.Range("I136").Value = 3831.2
$ws.Range("K136").Value = 11493.6
$ws.Range("M136").Value = -8943.599999999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 28622.688
$ws.Range("I20").Value = 36755.668
$ws.Range("K20").Value = 36755.668
$ws.Range("M20").Value = -36508.668
$ws.Range("H82").Value = 34281
$ws.Range("I82").Value = 21193.625
$ws.Range("J82").Value = 69180.664
$ws.Range("K82").Value = 21193.625
$ws.Range("L82").Value = 69180.664
$ws.Range("M82").Value = -20810.625
$ws.Range("N82").Value = -69946.664
$ws.Range("H85").Value = 34281
$ws.Range("I85").Value = 21193.625
$ws.Range("J85").Value = 69180.664
$ws.Range("K85").Value = 21193.625
$ws.Range("L85").Value = 69180.664
$ws.Range("M85").Value = -19867.625
$ws.Range("N85").Value = -71832.664
$ws.Range("H134").Value = 10941893
$ws.Range("I134").Value = 2383126.8
$ws.Range("K134").Value = 7149380.399999999
$ws.Range("M134").Value = -7146845.399999999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2004.3334
$ws.Range("I22").Value = 2205.2
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 2205.2
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -1855.2
$ws.Range("N22").Value = -1700
$ws.Range("H31").Value = 6667.1343
$ws.Range("I31").Value = 3485.818
$ws.Range("J31").Value = 10350.763
$ws.Range("K31").Value = 3485.818
$ws.Range("L31").Value = 10350.763
$ws.Range("M31").Value = -3190.818
$ws.Range("N31").Value = -10940.763
$ws.Range("H34").Value = 6667.1343
$ws.Range("I34").Value = 3485.818
$ws.Range("J34").Value = 10350.763
$ws.Range("K34").Value = 3485.818
$ws.Range("L34").Value = 10350.763
$ws.Range("M34").Value = -3283.818
$ws.Range("N34").Value = -10754.763
$ws.Range("H58").Value = 2231.6843
$ws.Range("I58").Value = 1809.3636
$ws.Range("K58").Value = 1809.3636
$ws.Range("M58").Value = -1606.3636
$ws.Range("H80").Value = 94949
$ws.Range("J80").Value = 94949
$ws.Range("L80").Value = 94949
$ws.Range("N80").Value = -97195
$ws.Range("H83").Value = 94949
$ws.Range("J83").Value = 94949
$ws.Range("L83").Value = 284847
$ws.Range("N83").Value = -296079
$ws.Range("H94").Value = 1211.5714
$ws.Range("I94").Value = 1193.3
$ws.Range("J94").Value = 1257.25
$ws.Range("K94").Value = 1193.3
$ws.Range("L94").Value = 1257.25
$ws.Range("M94").Value = -742.3
$ws.Range("N94").Value = -2159.25
$ws.Range("H124").Value = 66871.86
$ws.Range("J124").Value = 66871.86
$ws.Range("L124").Value = 66871.86
$ws.Range("N124").Value = -71781.86
$ws.Range("H125").Value = 63984
$ws.Range("J125").Value = 63984
$ws.Range("L125").Value = 63984
$ws.Range("N125").Value = -68904
$ws.Range("H132").Value = 3181.1875
$ws.Range("I132").Value = 2000.3846
$ws.Range("K132").Value = 6001.1538
$ws.Range("M132").Value = -3471.1538
$ws.Range("H134").Value = 2098.1924
$ws.Range("I134").Value = 1607.0238
$ws.Range("K134").Value = 4821.0714
$ws.Range("M134").Value = -2286.0714
$ws.Range("H136").Value = 2231.6843
$ws.Range("I136").Value = 1809.3636
$ws.Range("K136").Value = 5428.0908
$ws.Range("M136").Value = -2878.0908

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 10558.637
$ws.Range("I14").Value = 10558.637
$ws.Range("K14").Value = 31675.911
$ws.Range("M14").Value = -31502.911
$ws.Range("H75").Value = 729
$ws.Range("J75").Value = 701.7143
$ws.Range("L75").Value = 2105.1429
$ws.Range("N75").Value = -4101.1429
$ws.Range("H78").Value = 729
$ws.Range("J78").Value = 701.7143
$ws.Range("L78").Value = 6315.428699999999
$ws.Range("N78").Value = -16299.4287
$ws.Range("H113").Value = 4097.3125
$ws.Range("I113").Value = 4299.3335
$ws.Range("J113").Value = 4050.6924
$ws.Range("K113").Value = 12898.0005
$ws.Range("L113").Value = 12152.0772
$ws.Range("M113").Value = -10728.0005
$ws.Range("N113").Value = -16492.0772
$ws.Range("H131").Value = 5799.8
$ws.Range("J131").Value = 4749.75
$ws.Range("L131").Value = 14249.25
$ws.Range("N131").Value = -24329.25

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 156.75
$ws.Range("I107").Value = 176.83333
$ws.Range("J107").Value = 96.5
$ws.Range("K107").Value = 176.83333
$ws.Range("L107").Value = 96.5
$ws.Range("M107").Value = 1743.16667
$ws.Range("N107").Value = -3936.5
$ws.Range("H116").Value = 104330
$ws.Range("J116").Value = 104995
$ws.Range("L116").Value = 104995
$ws.Range("N116").Value = -114173
$ws.Range("H134").Value = 88471.5
$ws.Range("J134").Value = 88471.5
$ws.Range("L134").Value = 265414.5
$ws.Range("N134").Value = -270484.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2531.5
$ws.Range("I100").Value = 2531.5
$ws.Range("K100").Value = 2531.5
$ws.Range("M100").Value = -1990.5
$ws.Range("H122").Value = 1920202.2
$ws.Range("I122").Value = 3835405
$ws.Range("K122").Value = 11506215
$ws.Range("M122").Value = -11503765
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H137").Value = 98886.75
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 98886.75
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 98886.75
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -109086.75

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 8261.5
$ws.Range("J45").Value = 6848.6665
$ws.Range("L45").Value = 6848.6665
$ws.Range("N45").Value = -7830.6665
$ws.Range("H93").Value = 75891.5
$ws.Range("J93").Value = 75891.5
$ws.Range("L93").Value = 75891.5
$ws.Range("N93").Value = -80883.5
$ws.Range("H97").Value = 27475
$ws.Range("J97").Value = 27475
$ws.Range("L97").Value = 27475
$ws.Range("N97").Value = -29457
$ws.Range("H108").Value = 114000
$ws.Range("J108").Value = 114000
$ws.Range("L108").Value = 114000
$ws.Range("N108").Value = -121680
$ws.Range("H114").Value = 52989.332
$ws.Range("J114").Value = 52989.332
$ws.Range("L114").Value = 52989.332
$ws.Range("N114").Value = -61667.332
$ws.Range("H124").Value = 65381.75
$ws.Range("J124").Value = 65381.75
$ws.Range("L124").Value = 65381.75
$ws.Range("N124").Value = -75201.75
$ws.Range("H135").Value = 96249.25
$ws.Range("J135").Value = 96249.25
$ws.Range("L135").Value = 96249.25
$ws.Range("N135").Value = -106389.25

Write-Host "Applied 240 cell updates and 2 clears."
